# Natmi following Dr Hou advice
# Rebuild the LR-pairs sheet (Pf4 -> Cxcr3) to include all sending/target
# cluster combinations (ECs, M2, sCs) instead of just the single M2->M2 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: col A = Sending cluster, D = Target cluster (B/C ligand/receptor
# symbol are always Pf4 / Cxcr3 for this sheet).
$rows = @(
    @("ECs", "Pf4", "Cxcr3", "ECs", 2, 0.6666666666666666, 26.70522333333333, 80.11566999999999, 0.2696387810903721, 0.2696387810903721, 1, 0.3333333333333333, 1.863797, 5.591391, 0.5455000708290748, 0.5455000708290748, 49.77311513299667, 447.9580361969699, 0.1470879741830634, 0.1470879741830634),
    @("ECs", "Pf4", "Cxcr3", "M2",  2, 0.6666666666666666, 26.70522333333333, 80.11566999999999, 0.2696387810903721, 0.2696387810903721, 3, 1, 1.552879, 4.658637000000001, 0.4544999291709252, 0.4544999291709252, 41.46998050464334, 373.22982454179, 0.1225508069073087, 0.1225508069073087),
    @("M2",  "Pf4", "Cxcr3", "ECs", 3, 1, 72.30024833333333, 216.900745, 0.7300051600316594, 0.7300051600316595, 1, 0.3333333333333333, 1.863797, 5.591391, 0.5455000708290748, 0.5455000708290748, 134.7529859429216, 1212.776873486295, 0.3982178665028603, 0.3982178665028603),
    @("M2",  "Pf4", "Cxcr3", "M2",  3, 1, 72.30024833333333, 216.900745, 0.7300051600316594, 0.7300051600316595, 3, 1, 1.552879, 4.658637000000001, 0.4544999291709252, 0.4544999291709252, 112.2735373316183, 1010.461835984565, 0.3317872935287992, 0.3317872935287992),
    @("sCs", "Pf4", "Cxcr3", "ECs", 1, 0.3333333333333333, 0.03526433333333334, 0.105793, 0.000356058877968489, 0.0003560588779684891, 1, 0.3333333333333333, 1.863797, 5.591391, 0.5455000708290748, 0.5455000708290748, 0.06572555867366667, 0.5915300280629999, 0.0001942301431511317, 0.0001942301431511317),
    @("sCs", "Pf4", "Cxcr3", "M2",  1, 0.3333333333333333, 0.03526433333333334, 0.105793, 0.000356058877968489, 0.0003560588779684891, 3, 1, 1.552879, 4.658637000000001, 0.4544999291709252, 0.4544999291709252, 0.05476124268233334, 0.492851184141, 0.0001618287348173574, 0.0001618287348173574)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
